$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.434.66'
$ws.Range("E2").Value = '  +2.32%  '
$ws.Range("D3").Value = '2.066.87'
$ws.Range("E3").Value = '  +5.66%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'236.27"
$ws.Range("E5").Value = '  +2.49%  '
$ws.Range("D6").Value = "'0.617"
$ws.Range("E6").Value = '  +3.77%  '
$ws.Range("D7").Value = "'57.93"
$ws.Range("E7").Value = '  +9.65%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("E9").Value = '  +4.87%  '
$ws.Range("E10").Value = '  +1.25%  '
$ws.Range("D11").Value = "'0.0761"
$ws.Range("E11").Value = '  +3.97%  '
$ws.Range("E12").Value = '  +4.50%  '
$ws.Range("D13").Value = '2.374.16'
$ws.Range("E13").Value = '  +5.73%  '
$ws.Range("D14").Value = "'14.32"
$ws.Range("E14").Value = '  +4.51%  '
$ws.Range("D15").Value = "'20.90"
$ws.Range("E15").Value = '  +6.83%  '
$ws.Range("D16").Value = "'0.777"
$ws.Range("E16").Value = '  +5.23%  '
$ws.Range("E17").Value = '  +4.98%  '
$ws.Range("D18").Value = '2.065.65'
$ws.Range("E18").Value = '  +5.49%  '
$ws.Range("D19").Value = '37.615.05'
$ws.Range("E19").Value = '  +2.97%  '
$ws.Range("D20").Value = "'6.10"
$ws.Range("E20").Value = '  +23.66%  '
$ws.Range("D21").Value = "'68.56"
$ws.Range("E21").Value = '  +2.30%  '
$ws.Range("D22").Value = '0.0₃0811'
$ws.Range("E22").Value = '  +3.10%  '
$ws.Range("D23").Value = "'224.81"
$ws.Range("E23").Value = '  +2.39%  '
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = '  -0.16%  '
$ws.Range("D25").Value = "'2.45"
$ws.Range("E25").Value = '  +6.81%  '
$ws.Range("E26").Value = '  +3.38%  '
$ws.Range("D27").Value = "'162.81"
$ws.Range("E27").Value = '  +1.75%  '
$ws.Range("E28").Value = '  +4.90%  '
$ws.Range("E29").Value = '  +8.51%  '
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").Value = "'1.39"
$ws.Range("E30").Value = '  +8.65%  '
$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").Value = "'19.29"
$ws.Range("E31").Value = '  +2.96%  '
$ws.Range("E32").Value = '  +2.97%  '
$ws.Range("D33").Value = "'2.62"
$ws.Range("E33").Value = '  +16.62%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = "'4.48"
$ws.Range("E34").Value = '  +4.52%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = "'0.0629"
$ws.Range("E35").Value = '  +5.61%  '
$ws.Range("E36").Value = '  +7.91%  '
$ws.Range("E37").Value = '  -0.06%  '
$ws.Range("E38").Value = '  +0.64%  '
$ws.Range("D39").Value = "'3.35"
$ws.Range("E39").Value = '  +5.98%  '
$ws.Range("E40").Value = '  +15.46%  '
$ws.Range("E41").Value = '  -1.97%  '
$ws.Range("D42").Value = "'4.47"
$ws.Range("E42").Value = '  +31.38%  '
$ws.Range("D43").Value = "'0.0953"
$ws.Range("E43").Value = '  +10.32%  '
$ws.Range("D44").Value = '1.471.09'
$ws.Range("E44").Value = '  +5.23%  '
$ws.Range("D45").Value = "'95.48"
$ws.Range("E45").Value = '  +11.15%  '
$ws.Range("E46").Value = '  +5.88%  '
$ws.Range("D47").Value = "'16.16"
$ws.Range("E47").Value = '  +10.34%  '
$ws.Range("E48").Value = '  +5.64%  '
$ws.Range("E49").Value = '  +9.75%  '
$ws.Range("E50").Value = '  +4.44%  '
$ws.Range("E51").Value = '  +2.66%  '

# Reset style index on cells forced to text so no stray style is introduced
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
